$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 160.8
$ws.Range("I8").Value = 191
$ws.Range("K8").Value = 573
$ws.Range("M8").Value = -434
$ws.Range("H38").Value = 2110.3076
$ws.Range("I38").Value = 304.25
$ws.Range("K38").Value = 912.75
$ws.Range("M38").Value = -540.75
$ws.Range("H55").Value = 281.33334
$ws.Range("I55").Value = 281.33334
$ws.Range("K55").Value = 281.33334
$ws.Range("M55").Value = -67.33334000000002
$ws.Range("H74").Value = 8194.962
$ws.Range("I74").Value = 6676.643
$ws.Range("J74").Value = 9966.333
$ws.Range("K74").Value = 6676.643
$ws.Range("L74").Value = 9966.333
$ws.Range("M74").Value = -5740.643
$ws.Range("N74").Value = -11838.333
$ws.Range("H77").Value = 8194.962
$ws.Range("I77").Value = 6676.643
$ws.Range("J77").Value = 9966.333
$ws.Range("K77").Value = 33383.215
$ws.Range("L77").Value = 49831.665
$ws.Range("M77").Value = -28703.215
$ws.Range("N77").Value = -59191.665
$ws.Range("H127").Value = 1782
$ws.Range("I127").Value = 991.3333
$ws.Range("J127").Value = 2375
$ws.Range("K127").Value = 2973.9999
$ws.Range("L127").Value = 7125
$ws.Range("M127").Value = 1986.0001
$ws.Range("N127").Value = -17045
$ws.Range("H129").Value = 2855.0588
$ws.Range("I129").Value = 2709.1538
$ws.Range("K129").Value = 8127.4614
$ws.Range("M129").Value = -3127.4614
$ws.Range("H132").Value = 322425.4
$ws.Range("I132").Value = 346727.28
$ws.Range("J132").Value = 6501.25
$ws.Range("K132").Value = 1040181.84
$ws.Range("L132").Value = 19503.75
$ws.Range("M132").Value = -1037651.84
$ws.Range("N132").Value = -24563.75
$ws.Range("H137").Value = 3545.577
$ws.Range("I137").Value = 2524.7144
$ws.Range("J137").Value = 3921.6843
$ws.Range("K137").Value = 7574.1432
$ws.Range("L137").Value = 11765.0529
$ws.Range("M137").Value = -5024.1432
$ws.Range("N137").Value = -16865.0529
$ws.Range("H138").Value = 3931.6086
$ws.Range("I138").Value = 3428.889
$ws.Range("J138").Value = 4053.8918
$ws.Range("K138").Value = 10286.667
$ws.Range("L138").Value = 12161.6754
$ws.Range("M138").Value = -5146.667000000001
$ws.Range("N138").Value = -22441.6754
$ws.Range("H140").Value = 79497
$ws.Range("J140").Value = 79497
$ws.Range("L140").Value = 79497
$ws.Range("N140").Value = -89857
$ws.Range("H141").Value = 1695.3334
$ws.Range("I141").Value = 1260.5555
$ws.Range("J141").Value = 2999.6667
$ws.Range("K141").Value = 3781.6665
$ws.Range("L141").Value = 8999.000100000001
$ws.Range("M141").Value = 1398.3335
$ws.Range("N141").Value = -19359.0001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8937335
$ws.Range("I32").Value = 10423119
$ws.Range("K32").Value = 10423119
$ws.Range("M32").Value = -10422832
$ws.Range("H88").Value = 3009.7273
$ws.Range("J88").Value = 3440.3333
$ws.Range("L88").Value = 3440.3333
$ws.Range("N88").Value = -4252.3333
$ws.Range("H91").Value = 3009.7273
$ws.Range("J91").Value = 3440.3333
$ws.Range("L91").Value = 3440.3333
$ws.Range("N91").Value = -6248.3333
$ws.Range("H122").Value = 2971.3
$ws.Range("I122").Value = 2464.125
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 7392.375
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -4942.375
$ws.Range("N122").Value = -19900
$ws.Range("H132").Value = 755592.2
$ws.Range("I132").Value = 1111518.8
$ws.Range("K132").Value = 3334556.4
$ws.Range("M132").Value = -3332026.4

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2064.0908
$ws.Range("I86").Value = 2116.8572
$ws.Range("J86").Value = 1971.75
$ws.Range("K86").Value = 2116.8572
$ws.Range("L86").Value = 1971.75
$ws.Range("M86").Value = -993.8571999999999
$ws.Range("N86").Value = -4217.75
$ws.Range("H89").Value = 2064.0908
$ws.Range("I89").Value = 2116.8572
$ws.Range("J89").Value = 1971.75
$ws.Range("K89").Value = 10584.286
$ws.Range("L89").Value = 9858.75
$ws.Range("M89").Value = -4968.286
$ws.Range("N89").Value = -21090.75
$ws.Range("H105").Value = 3821.5557
$ws.Range("I105").Value = 3347.1428
$ws.Range("K105").Value = 3347.1428
$ws.Range("M105").Value = -1600.1428

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1331.6428
$ws.Range("I16").Value = 1339.3
$ws.Range("J16").Value = 1312.5
$ws.Range("K16").Value = 1339.3
$ws.Range("L16").Value = 1312.5
$ws.Range("M16").Value = -1052.3
$ws.Range("N16").Value = -1886.5
$ws.Range("H58").Value = 826401.94
$ws.Range("I58").Value = 1236203
$ws.Range("J58").Value = 6799.8
$ws.Range("K58").Value = 1236203
$ws.Range("L58").Value = 6799.8
$ws.Range("M58").Value = -1236000
$ws.Range("N58").Value = -7205.8
$ws.Range("H62").Value = 4813.5713
$ws.Range("I62").Value = 4674.5
$ws.Range("J62").Value = 4999
$ws.Range("K62").Value = 4674.5
$ws.Range("L62").Value = 4999
$ws.Range("M62").Value = -4050.5
$ws.Range("N62").Value = -6247
$ws.Range("H65").Value = 4813.5713
$ws.Range("I65").Value = 4674.5
$ws.Range("J65").Value = 4999
$ws.Range("K65").Value = 23372.5
$ws.Range("L65").Value = 24995
$ws.Range("M65").Value = -20252.5
$ws.Range("N65").Value = -31235
$ws.Range("H94").Value = 18547.084
$ws.Range("I94").Value = 41044.2
$ws.Range("K94").Value = 41044.2
$ws.Range("M94").Value = -40593.2
$ws.Range("H105").Value = 51517.715
$ws.Range("I105").Value = 59707.332
$ws.Range("J105").Value = 2380
$ws.Range("K105").Value = 59707.332
$ws.Range("L105").Value = 2380
$ws.Range("M105").Value = -57960.332
$ws.Range("N105").Value = -5874
$ws.Range("H107").Value = 1955.6471
$ws.Range("J107").Value = 2990.5
$ws.Range("L107").Value = 2990.5
$ws.Range("N107").Value = -6830.5
$ws.Range("H113").Value = 1331.6428
$ws.Range("I113").Value = 1339.3
$ws.Range("J113").Value = 1312.5
$ws.Range("K113").Value = 1339.3
$ws.Range("L113").Value = 1312.5
$ws.Range("M113").Value = 830.7
$ws.Range("N113").Value = -5652.5
$ws.Range("H122").Value = 2097.8572
$ws.Range("I122").Value = 1907.5
$ws.Range("J122").Value = 2478.5715
$ws.Range("K122").Value = 5722.5
$ws.Range("L122").Value = 7435.7145
$ws.Range("M122").Value = -3272.5
$ws.Range("N122").Value = -12335.7145
$ws.Range("H132").Value = 9632400
$ws.Range("I132").Value = 23646.188
$ws.Range("J132").Value = 25006408
$ws.Range("K132").Value = 70938.564
$ws.Range("L132").Value = 75019224
$ws.Range("M132").Value = -68408.564
$ws.Range("N132").Value = -75024284
$ws.Range("H134").Value = 3878.1667
$ws.Range("I134").Value = 1394.6086
$ws.Range("J134").Value = 61000
$ws.Range("K134").Value = 4183.825800000001
$ws.Range("L134").Value = 183000
$ws.Range("M134").Value = -1648.825800000001
$ws.Range("N134").Value = -188070
$ws.Range("H136").Value = 826401.94
$ws.Range("I136").Value = 1236203
$ws.Range("J136").Value = 6799.8
$ws.Range("K136").Value = 3708609
$ws.Range("L136").Value = 20399.4
$ws.Range("M136").Value = -3706059
$ws.Range("N136").Value = -25499.4

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 36.166668
$ws.Range("J2").Value = 32.714287
$ws.Range("L2").Value = 196.285722
$ws.Range("N2").Value = -422.285722
$ws.Range("H96").Value = 7322.4
$ws.Range("J96").Value = 8028
$ws.Range("L96").Value = 24084
$ws.Range("N96").Value = -28202

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2959.818
$ws.Range("J97").Value = 3055.8
$ws.Range("L97").Value = 3055.8
$ws.Range("N97").Value = -4047.8

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4194.5
$ws.Range("I7").Value = 3960.5557
$ws.Range("J7").Value = 5247.25
$ws.Range("K7").Value = 3960.5557
$ws.Range("L7").Value = 5247.25
$ws.Range("M7").Value = -3848.5557
$ws.Range("N7").Value = -5471.25
$ws.Range("H46").Value = 3009.3333
$ws.Range("I46").Value = 1257.1666
$ws.Range("K46").Value = 1257.1666
$ws.Range("M46").Value = -1069.1666
$ws.Range("H55").Value = 1564.6666
$ws.Range("I55").Value = 622.2
$ws.Range("J55").Value = 2035.9
$ws.Range("K55").Value = 622.2
$ws.Range("L55").Value = 2035.9
$ws.Range("M55").Value = -449.2
$ws.Range("N55").Value = -2381.9
$ws.Range("H100").Value = 11343.357
$ws.Range("J100").Value = 14800.4
$ws.Range("L100").Value = 14800.4
$ws.Range("N100").Value = -15882.4
$ws.Range("H122").Value = 3486.1482
$ws.Range("I122").Value = 3353.1667
$ws.Range("K122").Value = 10059.5001
$ws.Range("M122").Value = -7609.500100000001
$ws.Range("H126").Value = 4194.5
$ws.Range("I126").Value = 3960.5557
$ws.Range("J126").Value = 5247.25
$ws.Range("K126").Value = 11881.6671
$ws.Range("L126").Value = 15741.75
$ws.Range("M126").Value = -9411.667099999999
$ws.Range("N126").Value = -20681.75

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 2921.3635
$ws.Range("I100").Value = 3223.875
$ws.Range("J100").Value = 2114.6667
$ws.Range("K100").Value = 6447.75
$ws.Range("L100").Value = 4229.3334
$ws.Range("M100").Value = -5906.75
$ws.Range("N100").Value = -5311.3334
$ws.Range("H136").Value = 29301796
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

